# Auto-generated Excel COM-interop script
# Updates 'F' (想去人数 / want-to-go count) and a couple of G/D cells
# across the four sheets, per the upstream data refresh commit.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 40
$ws.Range("F4").Value = 8034
$ws.Range("F7").Value = 78
$ws.Range("F8").Value = 6947
$ws.Range("F9").Value = 6947
$ws.Range("F10").Value = 1131
$ws.Range("F11").Value = 519
$ws.Range("F14").Value = 692
$ws.Range("F16").Value = 303
$ws.Range("F17").Value = 160
$ws.Range("F18").Value = 226
$ws.Range("F19").Value = 167
$ws.Range("F20").Value = 119
$ws.Range("F21").Value = 11346
$ws.Range("F22").Value = 97
$ws.Range("F23").Value = 2167
$ws.Range("F25").Value = 3002
$ws.Range("F26").Value = 50
$ws.Range("F28").Value = 2603
$ws.Range("F31").Value = 266
$ws.Range("F32").Value = 40
$ws.Range("F34").Value = 2327
$ws.Range("F35").Value = 336
$ws.Range("F36").Value = 1581
$ws.Range("D37").Value = "北苑东路88号近中铁国际城小区 蓝地时尚庄园"
$ws.Range("G37").Value = 20
$ws.Range("F38").Value = 80
$ws.Range("F39").Value = 5714
$ws.Range("F40").Value = 1757
$ws.Range("F42").Value = 818
$ws.Range("F43").Value = 154
$ws.Range("F45").Value = 1117
$ws.Range("F47").Value = 1058
$ws.Range("F48").Value = 1495
$ws.Range("F49").Value = 94

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 5
$ws.Range("F8").Value = 248
$ws.Range("F23").Value = 5

# --- Sheet: 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 194
$ws.Range("F3").Value = 317

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 40
$ws.Range("F4").Value = 194
$ws.Range("F5").Value = 317
$ws.Range("F7").Value = 5
$ws.Range("F8").Value = 8034
$ws.Range("F11").Value = 78
$ws.Range("F12").Value = 6948
$ws.Range("F13").Value = 1131
$ws.Range("F14").Value = 519
$ws.Range("F16").Value = 692
$ws.Range("F18").Value = 303
$ws.Range("F19").Value = 160
$ws.Range("F20").Value = 226
$ws.Range("F21").Value = 167
$ws.Range("F24").Value = 11346
$ws.Range("F25").Value = 97
$ws.Range("F26").Value = 2167
$ws.Range("F27").Value = 2167
$ws.Range("F28").Value = 3002
$ws.Range("F29").Value = 2603
$ws.Range("F31").Value = 266
$ws.Range("F32").Value = 40
$ws.Range("F34").Value = 2327
$ws.Range("F35").Value = 336
$ws.Range("F36").Value = 1581
$ws.Range("D37").Value = "北苑东路88号近中铁国际城小区 蓝地时尚庄园"
$ws.Range("F37").Value = 69
$ws.Range("G37").Value = 20
$ws.Range("F38").Value = 80
$ws.Range("F39").Value = 5715
$ws.Range("F41").Value = 1757
$ws.Range("F44").Value = 818
$ws.Range("F45").Value = 154
$ws.Range("F46").Value = 1117
$ws.Range("F48").Value = 1058
$ws.Range("F49").Value = 1495
$ws.Range("F50").Value = 94
